$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 242.245169555222
$ws.Range("B2").Value = 159.7870462081287
$ws.Range("C2").Value = 1820.392602918578
$ws.Range("D2").Value = 231.3724111404355
$ws.Range("E2").Value = 192.8945719094686
